$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "results": insert a new "S*-unmerged" column (before the old
# "var" column) and a new "S*-unmergedND" column at the end.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("results")

# Insert a new column G; existing G:K shift right to H:L.
$ws1.Columns("G:G").Insert()
$ws1.Range("G1").Value = "S*-unmerged"
$ws1.Range("G2").Value = 194

# Append a new trailing column M, copying L's style first.
$ws1.Range("L1").Copy($ws1.Range("M1"))
$ws1.Range("M1").Value = "S*-unmergedND"
$ws1.Range("M2").Value = $false

# The old I2 (TRUE) is now at J2 after the column insert; the refreshed
# run shows it FALSE like all the other ND flags.
$ws1.Range("J2").Value = $false

# ---------------------------------------------------------------------
# Sheet "stats": insert a new "S*-unmerged" row before each "Kruskal"
# row, then refresh every numeric result with the latest benchmark run.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("stats")

# First block: new row 6 (pushes old row6 "Kruskal" down to row7).
# NB: the new row sits inside the A2:A6 merge, so its A/B cells must be
# restyled via BorderAround (not Range.Copy, which silently unmerges the
# destination) to pick up the same bold+bordered look as the rest of the
# column without breaking the merge.
$ws2.Rows("6:6").Insert()
$ws2.Range("A6").BorderAround(1)
$ws2.Range("B6").BorderAround(1)
$ws2.Range("B6").Value = "S*-unmerged"

# Second block: new row 12 (pushes old row11 "Kruskal" down to row13).
$ws2.Rows("12:12").Insert()
$ws2.Range("A12").BorderAround(1)
$ws2.Range("B12").BorderAround(1)
$ws2.Range("B12").Value = "S*-unmerged"

# Refresh block 1 (rows 2-7) with the new benchmark values.
$ws2.Range("C2").Value = 75
$ws2.Range("D2").Value = 0.0001129219308495522
$ws2.Range("E2").Value = 0.0268957051448524
$ws2.Range("F2").Value = 75
$ws2.Range("G2").Value = 0.003144462592899799
$ws2.Range("H2").Value = 0.004684991668909788
$ws2.Range("I2").Value = 0.005966922733932734
$ws2.Range("J2").Value = 0.009391450323164463
$ws2.Range("K2").Value = 0.0009838948026299477

$ws2.Range("C3").Value = 75
$ws2.Range("D3").Value = 0.002007377333939075
$ws2.Range("E3").Value = 0.03249291703104973
$ws2.Range("F3").Value = 75
$ws2.Range("G3").Value = 0.002872450277209282
$ws2.Range("H3").Value = 0.00690627982839942
$ws2.Range("I3").Value = 0.006613369099795818
$ws2.Range("J3").Value = 0.01225536083802581
$ws2.Range("K3").Value = 0.0010749576613307

$ws2.Range("C4").Value = 75
$ws2.Range("D4").Value = 0.00332916621118784
$ws2.Range("E4").Value = 0.04765735613182187
$ws2.Range("F4").Value = 75
$ws2.Range("G4").Value = 0.004393088165670633
$ws2.Range("H4").Value = 0.009997188113629818
$ws2.Range("I4").Value = 0.01035538129508495
$ws2.Range("J4").Value = 0.0170828877016902
$ws2.Range("K4").Value = 0.001668364740908146

$ws2.Range("C5").Value = 75
$ws2.Range("D5").Value = 0.0002576448023319244
$ws2.Range("E5").Value = 0.04365428956225514
$ws2.Range("F5").Value = 75
$ws2.Range("G5").Value = 0.004459449555724859
$ws2.Range("H5").Value = 0.007149436045438051
$ws2.Range("I5").Value = 0.01101410528644919
$ws2.Range("J5").Value = 0.01492588361725211
$ws2.Range("K5").Value = 0.00161285512149334

$ws2.Range("C6").Value = 161
$ws2.Range("D6").Value = 0.003033468965440989
$ws2.Range("E6").Value = 0.09924222994595766
$ws2.Range("F6").Value = 161
$ws2.Range("G6").Value = 0.00554067874327302
$ws2.Range("H6").Value = 0.01238089008256793
$ws2.Range("I6").Value = 0.05376495467498899
$ws2.Range("J6").Value = 0.0194350890815258
$ws2.Range("K6").Value = 0.002341007348150015

$ws2.Range("C7").Value = 1700
$ws2.Range("E7").Value = 0.01802199811208993

# Refresh block 2 (rows 8-13) -- identical data to block 1.
$ws2.Range("C8").Value = 75
$ws2.Range("D8").Value = 0.0001129219308495522
$ws2.Range("E8").Value = 0.0268957051448524
$ws2.Range("F8").Value = 75
$ws2.Range("G8").Value = 0.003144462592899799
$ws2.Range("H8").Value = 0.004684991668909788
$ws2.Range("I8").Value = 0.005966922733932734
$ws2.Range("J8").Value = 0.009391450323164463
$ws2.Range("K8").Value = 0.0009838948026299477

$ws2.Range("C9").Value = 75
$ws2.Range("D9").Value = 0.002007377333939075
$ws2.Range("E9").Value = 0.03249291703104973
$ws2.Range("F9").Value = 75
$ws2.Range("G9").Value = 0.002872450277209282
$ws2.Range("H9").Value = 0.00690627982839942
$ws2.Range("I9").Value = 0.006613369099795818
$ws2.Range("J9").Value = 0.01225536083802581
$ws2.Range("K9").Value = 0.0010749576613307

$ws2.Range("C10").Value = 75
$ws2.Range("D10").Value = 0.00332916621118784
$ws2.Range("E10").Value = 0.04765735613182187
$ws2.Range("F10").Value = 75
$ws2.Range("G10").Value = 0.004393088165670633
$ws2.Range("H10").Value = 0.009997188113629818
$ws2.Range("I10").Value = 0.01035538129508495
$ws2.Range("J10").Value = 0.0170828877016902
$ws2.Range("K10").Value = 0.001668364740908146

$ws2.Range("C11").Value = 75
$ws2.Range("D11").Value = 0.0002576448023319244
$ws2.Range("E11").Value = 0.04365428956225514
$ws2.Range("F11").Value = 75
$ws2.Range("G11").Value = 0.004459449555724859
$ws2.Range("H11").Value = 0.007149436045438051
$ws2.Range("I11").Value = 0.01101410528644919
$ws2.Range("J11").Value = 0.01492588361725211
$ws2.Range("K11").Value = 0.00161285512149334

$ws2.Range("C12").Value = 161
$ws2.Range("D12").Value = 0.003033468965440989
$ws2.Range("E12").Value = 0.09924222994595766
$ws2.Range("F12").Value = 161
$ws2.Range("G12").Value = 0.00554067874327302
$ws2.Range("H12").Value = 0.01238089008256793
$ws2.Range("I12").Value = 0.05376495467498899
$ws2.Range("J12").Value = 0.0194350890815258
$ws2.Range("K12").Value = 0.002341007348150015

$ws2.Range("C13").Value = 1700
$ws2.Range("E13").Value = 0.01802199811208993
